$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.8971247825105778
$ws.Range("D2").Value = 0.9197820816801586
$ws.Range("C3").Value = 0.8126266953362704
$ws.Range("D3").Value = 0.8512590892885307
$ws.Range("C4").Value = 0.7425063532679597
$ws.Range("D4").Value = 0.7926889049275229
$ws.Range("C5").Value = 0.6838220640019288
$ws.Range("D5").Value = 0.7422349058032773
$ws.Range("B6").Value = 0.6672664546811148
$ws.Range("C6").Value = 0.6343786487394392
$ws.Range("D6").Value = 0.6987908319406669
$ws.Range("C7").Value = 0.591995898535953
$ws.Range("D7").Value = 0.6605717329158735
$ws.Range("C8").Value = 0.5558047588733271
$ws.Range("D8").Value = 0.6272914545507864
$ws.Range("C9").Value = 0.5247439588204366
$ws.Range("D9").Value = 0.5976139564099896
$ws.Range("C10").Value = 0.4983584500551801
$ws.Range("D10").Value = 0.5718147670247357
$ws.Range("C11").Value = 0.4750154045392584
$ws.Range("D11").Value = 0.5486653531852319
$ws.Range("B12").Value = 0.4925478298089613
$ws.Range("C12").Value = 0.455266362415196
$ws.Range("D12").Value = 0.5289753308491144
$ws.Range("B13").Value = 0.4746320322777658
$ws.Range("C13").Value = 0.4376891083665949
$ws.Range("D13").Value = 0.5115752281054643
$ws.Range("C14").Value = 0.4218876761366206
$ws.Range("D14").Value = 0.4955575564136658
$ws.Range("B15").Value = 0.4446738934837428
$ws.Range("C15").Value = 0.407647663132007
$ws.Range("D15").Value = 0.4812348569706433
$ws.Range("C16").Value = 0.3949563141524006
$ws.Range("D16").Value = 0.4683683986195139
$ws.Range("B17").Value = 0.4209194239897868
$ws.Range("C17").Value = 0.3837548360975112
$ws.Range("D17").Value = 0.4568443918027592
$ws.Range("C18").Value = 0.3733415513071138
$ws.Range("D18").Value = 0.4464853003180151
$ws.Range("C19").Value = 0.3640296964389892
$ws.Range("D19").Value = 0.4373592252574932
$ws.Range("B20").Value = 0.3938282540546766
$ws.Range("C20").Value = 0.3556796603725879
$ws.Range("D20").Value = 0.4287351722782844
$ws.Range("C21").Value = 0.3484394258424802
$ws.Range("D21").Value = 0.422046905593858
$ws.Range("C22").Value = 0.3420683858531907
$ws.Range("D22").Value = 0.4151194467820212
$ws.Range("C23").Value = 0.3367602058140226
$ws.Range("D23").Value = 0.4086381758744945
$ws.Range("C24").Value = 0.3316283972251199
$ws.Range("D24").Value = 0.4031783071516437
$ws.Range("B25").Value = 0.3639120098497534
$ws.Range("C25").Value = 0.327038342255182
$ws.Range("D25").Value = 0.3984626372962711
$ws.Range("C26").Value = 0.3229683687147928
$ws.Range("D26").Value = 0.3942600912225846
$ws.Range("B27").Value = 0.3555523258728327
$ws.Range("C27").Value = 0.3190783241702697
$ws.Range("D27").Value = 0.3903654773767701
$ws.Range("B28").Value = 0.3519424894369961
$ws.Range("C28").Value = 0.3152630257727591
$ws.Range("D28").Value = 0.3868533890160225
$ws.Range("B29").Value = 0.3486622624307191
$ws.Range("C29").Value = 0.3117989577803505
$ws.Range("D29").Value = 0.3838708480110637
$ws.Range("B30").Value = 0.3456791449892515
$ws.Range("C30").Value = 0.3084314110448401
$ws.Range("D30").Value = 0.3805749489510276
$ws.Range("C31").Value = 0.3060278216679467
$ws.Range("D31").Value = 0.3775913964713172
$ws.Range("C32").Value = 0.3036906888518494
$ws.Range("D32").Value = 0.3750793781179427
$ws.Range("C33").Value = 0.3016353999090092
$ws.Range("D33").Value = 0.3729953113507635
$ws.Range("C34").Value = 0.299777555256703
$ws.Range("D34").Value = 0.3711044596699966
$ws.Range("C35").Value = 0.2979489511249988
$ws.Range("D35").Value = 0.369393536032611
$ws.Range("C36").Value = 0.2962270095864155
$ws.Range("D36").Value = 0.3678664470830053
$ws.Range("B37").Value = 0.3310371625304339
$ws.Range("C37").Value = 0.2946692058340588
$ws.Range("D37").Value = 0.3665558917601686
$ws.Range("C38").Value = 0.2931037554667104
$ws.Range("D38").Value = 0.3651710998551696
$ws.Range("C39").Value = 0.2919079028975835
$ws.Range("D39").Value = 0.3639333536441638
$ws.Range("C40").Value = 0.2908248030583327
$ws.Range("D40").Value = 0.3628608054939622
$ws.Range("B41").Value = 0.3260242752213208
$ws.Range("C41").Value = 0.2897779876731053
$ws.Range("D41").Value = 0.3618832239208111
$ws.Range("C42").Value = 0.2887188566405427
$ws.Range("D42").Value = 0.3609919837001631
$ws.Range("B43").Value = 0.3241150343351163
$ws.Range("C43").Value = 0.2877521212731215
$ws.Range("D43").Value = 0.360178328517107
$ws.Range("C44").Value = 0.276538956931996
$ws.Range("D44").Value = 0.3517252388298802
